$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the existing data (rows 2-21, columns A-C) before shifting it down.
$oldData = $ws.Range("A2:C21").Value2

# Shift the existing 20 rows of data down by 9 rows -> new rows 11-30.
$ws.Range("A11:C30").Value2 = $oldData

# New rows 2-10 hold freshly added data (9 rows).
$topData = @(
    @(0.0609338097274303, 0.0045814891345798, 0.0629191175103187),
    @(0.0303905457258224, 0.0102319931611418, 0.0383317954838275),
    @(-0.0080939643085002, 0.0009162978967650999, -0.0238237436860799),
    @(0.0160352122038602, -0.0207694191485643, 0.0113010071218013),
    @(-0.0288633834570646, -0.00137444678694, -0.0123700210824608),
    @(0.0059559359215199, 0.0401643887162208, 0.0155770638957619),
    @(0.0117591563612222, 0.0131336031481623, 0.0174096599221229),
    @(0.0245873257517814, 0.0372627787292003, 0.028557950630784),
    @(0.0038179077673703, 0.0119118718430399, -0.0335975885391235)
)
$topArr = New-Object 'object[,]' 9,3
for ($i = 0; $i -lt 9; $i++) {
    for ($j = 0; $j -lt 3; $j++) {
        $topArr[$i, $j] = $topData[$i][$j]
    }
}
$ws.Range("A2:C10").Value2 = $topArr

# New row 31 holds one more freshly appended row of data.
$ws.Range("A31").Value2 = 0.0256563406437635
$ws.Range("B31").Value2 = 0.0372627787292003
$ws.Range("C31").Value2 = 0.0050396383740007
